$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.46%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.27%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.143"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.53%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'-0.70%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'FTXToken"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'1.680"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.21%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'MXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'0.9356"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.33%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'0.1201"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.73%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'WazirX"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.1819"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.24%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'MandalaExchangeToken"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.09040"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.99%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'BitrueCoin"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.04142"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.08%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'BitMartToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.1054"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.27%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitForexToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.001283"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.87%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'TigerCash"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.005830"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.91%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'LEO"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'3.341"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.15%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'GateToken"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'4.321"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.06%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.82%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.61%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.618"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.19%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1342"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.68%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2837"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.49%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-1.29%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001282"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.47%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-7.10%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'6.11%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02417"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-2.40%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05157"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.28%"
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'-1.54%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.92%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'15.46%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.003302"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'72.47%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007572"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.20%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3304"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.43%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006816"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.48%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.2731"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-35.59%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'35.25%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = "Normal"
